# Fruta / hortaliza, semanal
#
# The underlying data rows (2-26, 28-45; row 27 is untouched) were
# reshuffled: each row keeps its "identity" columns (Mercado, Región,
# Categoría, Variedad, Calidad, Unidad de comercialización, Origen,
# Kg o Unidades, Clasificación, ...) but is now paired with a different
# set of "observation" values - Fecha (D), Volumen (J), Precio mínimo
# (K), Precio máximo (L), Precio promedio ponderado (M) and Precio $/Kg
# (P) - taken from another row of the same original table.
#
# Each triple (target row -> source row) below describes which row's
# D/J/K/L/M/P values now land on the target row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @(2, 44181, 250, 1400, 1500, 1450, 725),
    @(3, 44235, 250, 4500, 5000, 4750, 2375),
    @(4, 44305, 300, 900, 1000, 950, 475),
    @(5, 44326, 300, 1400, 1500, 1450, 725),
    @(6, 44428, 270, 3500, 3800, 3650, 1825),
    @(7, 44424, 300, 2500, 3000, 2750, 1375),
    @(8, 44349, 300, 1800, 2000, 1900, 950),
    @(9, 44258, 150, 2400, 2500, 2450, 1225),
    @(10, 44431, 300, 1900, 2000, 1950, 975),
    @(11, 44435, 500, 1800, 2000, 1930, 965),
    @(12, 44169, 300, 2000, 2500, 2250, 1125),
    @(13, 44263, 270, 1900, 2000, 1950, 975),
    @(14, 44272, 250, 2800, 3000, 2900, 1450),
    @(15, 44298, 300, 1400, 1500, 1450, 725),
    @(16, 44319, 300, 1900, 2000, 1950, 975),
    @(17, 44293, 250, 1500, 1800, 1650, 825),
    @(18, 44320, 250, 1400, 1500, 1450, 725),
    @(19, 44284, 300, 1800, 2000, 1900, 950),
    @(20, 44417, 300, 3000, 3500, 3250, 1625),
    @(21, 44253, 300, 2400, 2500, 2450, 1225),
    @(22, 44334, 200, 2800, 3000, 2900, 1450),
    @(23, 44452, 300, 1900, 2000, 1950, 975),
    @(24, 44405, 300, 3800, 4000, 3900, 1950),
    @(25, 44442, 200, 2400, 2500, 2450, 1225),
    @(26, 44448, 270, 1900, 2000, 1950, 975),
    @(28, 44237, 200, 2500, 3000, 2750, 1375),
    @(29, 44266, 300, 1800, 2000, 1900, 950),
    @(30, 44165, 300, 1000, 1200, 1100, 550),
    @(31, 44249, 300, 2400, 2500, 2450, 1225),
    @(32, 44343, 300, 1500, 2000, 1750, 875),
    @(33, 44267, 300, 1400, 1500, 1450, 725),
    @(34, 44312, 300, 1000, 1200, 1100, 550),
    @(35, 44386, 250, 3500, 4000, 3750, 1875),
    @(36, 44243, 200, 2900, 3000, 2950, 1475),
    @(37, 44433, 200, 1800, 2000, 1900, 950),
    @(38, 44221, 200, 2900, 3000, 2950, 1475),
    @(39, 44279, 200, 1700, 1800, 1750, 875),
    @(40, 44397, 300, 3500, 4000, 3750, 1875),
    @(41, 44356, 300, 2400, 2500, 2450, 1225),
    @(42, 44203, 300, 2000, 2500, 2250, 1125),
    @(43, 44410, 250, 2800, 3000, 2900, 1450),
    @(44, 44176, 300, 1900, 2000, 1950, 975),
    @(45, 44323, 200, 2400, 2500, 2450, 1225)
)

foreach ($entry in $rowData) {
    $row = $entry[0]
    $ws.Range("D$row").Value = $entry[1]
    $ws.Range("J$row").Value = $entry[2]
    $ws.Range("K$row").Value = $entry[3]
    $ws.Range("L$row").Value = $entry[4]
    $ws.Range("M$row").Value = $entry[5]
    $ws.Range("P$row").Value = $entry[6]
}
